# Apply the "Reservas" sheet update: rows 2-19 get new/edited booking data.
# Columns: A=Fecha (date serial), B=Hora inicio, C=Hora fin (time fractions),
#          D=Profesor, E=Curso, F=Recurso, G=Observaciones

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reservas")

$rows = @(
    @{ r = 2;  a = 45845; b = 0.40625;            c = 0.46875;            d = "JOSEFINA ISABEL ESPINOSA BERRIOS";   e = "1° MEDIO A";  f = "ENLACE MEDIA";  g = "" },
    @{ r = 3;  a = 45845; b = 0.4791666666666667; c = 0.5416666666666666; d = "PILAR ALEJANDRA GONZÁLEZ OLIVARES";  e = "1° BÁSICO A"; f = "ENLACE MEDIA";  g = "letrapps" },
    @{ r = 4;  a = 45846; b = 0.3333333333333333; c = 0.3958333333333333; d = "MAGALY GIOCONDA MOLINA NÚÑEZ";       e = "3° BÁSICO A"; f = "ENLACE BASICA"; g = "umaximo" },
    @{ r = 5;  a = 45848; b = 0.3333333333333333; c = 0.3958333333333333; d = "JENNIFER SAMANTHA GONZÁLEZ RIVERA";  e = "8° BÁSICO B"; f = "ENLACE BASICA"; g = "cs. nat" },
    @{ r = 6;  a = 45848; b = 0.3333333333333333; c = 0.3958333333333333; d = "LETICIA PAOLA LÓPEZ IGLESIAS";       e = "5° BÁSICO A"; f = "ENLACE MEDIA";  g = "umaximo" },
    @{ r = 7;  a = 45848; b = 0.40625;            c = 0.46875;            d = "PILAR ALEJANDRA GONZÁLEZ OLIVARES";  e = "1° BÁSICO A"; f = "ENLACE MEDIA";  g = "umaximo" },
    @{ r = 8;  a = 45848; b = 0.4791666666666667; c = 0.5416666666666666; d = "PILAR ALEJANDRA GONZÁLEZ OLIVARES";  e = "6° BÁSICO A"; f = "ENLACE BASICA"; g = "umaximo" },
    @{ r = 9;  a = 45848; b = 0.5833333333333334; c = 0.6458333333333334; d = "JENNIFER SAMANTHA GONZÁLEZ RIVERA";  e = "8° BÁSICO A"; f = "ENLACE MEDIA";  g = "Cs. Nat." },
    @{ r = 10; a = 45848; b = 0.5833333333333334; c = 0.6458333333333334; d = "PILAR ALEJANDRA GONZÁLEZ OLIVARES";  e = "6° BÁSICO B"; f = "ENLACE BASICA"; g = "umaximo" },
    @{ r = 11; a = 45848; b = 0.5833333333333334; c = 0.6458333333333334; d = "MARTA YOLANDA PONCE LÓPEZ";          e = "5° BÁSICO B"; f = "TABLETS";       g = "" },
    @{ r = 12; a = 45849; b = 0.3333333333333333; c = 0.3958333333333333; d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE BASICA"; g = "" },
    @{ r = 13; a = 45845; b = 0.3333333333333333; c = 0.3958333333333333; d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "LAPICES 3D";    g = "" },
    @{ r = 14; a = 45847; b = 0.5833333333333334; c = 0.6458333333333334; d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE MEDIA";  g = "" },
    @{ r = 15; a = 45847; b = 0.65625;            c = 0.6875;             d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE MEDIA";  g = "" },
    @{ r = 16; a = 45847; b = 0.6875;             c = 0.7708333333333334; d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE MEDIA";  g = "" },
    @{ r = 17; a = 45847; b = 0.3333333333333333; c = 0.3958333333333333; d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE BASICA"; g = "" },
    @{ r = 18; a = 45847; b = 0.40625;            c = 0.46875;            d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE BASICA"; g = "" },
    @{ r = 19; a = 45847; b = 0.4791666666666667; c = 0.5416666666666666; d = "ALEJANDRA LORENA MUÑOZ TRUJILLO";    e = "1° BÁSICO A"; f = "ENLACE BASICA"; g = "" }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value = $row.a
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"

    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 2).NumberFormat = "h:mm:ss"

    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 3).NumberFormat = "h:mm:ss"

    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f

    # Always touch G — clears any stale leftover value from the old sheet
    # when the target observation is blank, writes the new text otherwise.
    if ($row.g -ne "") {
        $ws.Cells.Item($r, 7).Value = $row.g
    } else {
        $ws.Cells.Item($r, 7).ClearContents()
    }
}
